$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices for B,C,D,E,F,I,J,K,L,M
$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13)

# New vm_pu values per row (row 2..25), in the order of $cols above
$data = @{}
$data[2] = @(1.02, 1.045608419067008, 1.048080382523962, 1.043079053674487, 1.044447821620053, 1.044609206762441, 1.050667243542236, 1.050841225952017, 1.045853932001869, 1.047218844371849)
$data[3] = @(1.02, 1.047079355999303, 1.049255227248778, 1.044372084595887, 1.046546796370809, 1.045159702223752, 1.051783338939877, 1.051827104596299, 1.046956663991478, 1.049125703094741)
$data[4] = @(1.02, 1.048028127108215, 1.050012815318832, 1.045205566982232, 1.047901574026912, 1.045513029140961, 1.052502111585902, 1.052461852011942, 1.047666472505622, 1.050355808828114)
$data[5] = @(1.02, 1.048426278985594, 1.0503306884309, 1.045555208946008, 1.048470330578876, 1.045660883694322, 1.052803475951623, 1.052727946741115, 1.047963992580434, 1.050872065333887)
$data[6] = @(1.02, 1.048493089071075, 1.050384024812593, 1.045613871467583, 1.048565781495503, 1.045685669168332, 1.052854029284223, 1.052772581306626, 1.04801389600083, 1.050958696051607)
$data[7] = @(1.02, 1.048033450013649, 1.050017065167761, 1.045210241864725, 1.047909176864928, 1.045515007461984, 1.052506141592488, 1.052465410531197, 1.047670451442131, 1.05036271050581)
$data[8] = @(1.02, 1.046106162565492, 1.048477973780897, 1.043516707668984, 1.045157897029927, 1.044795849028816, 1.051045146544226, 1.051175072365079, 1.0462273845262, 1.047864066734267)
$data[9] = @(1.02, 1.04268632388365, 1.045745457720053, 1.040507522666601, 1.040282770532945, 1.043506290413033, 1.048444071644782, 1.048876568890983, 1.043655466099146, 1.043431440366485)
$data[10] = @(1.02, 1.040389706708623, 1.043909456219709, 1.038483925979756, 1.037013136474425, 1.042631234397535, 1.046691487209429, 1.047327031238536, 1.041920673191321, 1.040455119132723)
$data[11] = @(1.02, 1.039391103837639, 1.043110919184226, 1.037603386859923, 1.035592402469881, 1.042248605766153, 1.04592805898617, 1.046651860020344, 1.041164556553121, 1.039161017715103)
$data[12] = @(1.02, 1.039019539173341, 1.042813764316127, 1.037275654665207, 1.035063906518435, 1.042105914145073, 1.045643792029412, 1.046400427866246, 1.040882946347631, 1.038679504740738)
$data[13] = @(1.02, 1.039099270222074, 1.042877529784636, 1.037345984452274, 1.035177306065784, 1.04213654772035, 1.045704799941843, 1.046454390211755, 1.040943387054764, 1.038782828688911)
$data[14] = @(1.02, 1.039360403312709, 1.043086367406514, 1.037576309972281, 1.03554873273601, 1.042236822424938, 1.045904575656244, 1.046631089750366, 1.04114129403758, 1.039121232697872)
$data[15] = @(1.02, 1.039521210883129, 1.043214966907615, 1.037718133163988, 1.03577747785096, 1.042298529730845, 1.0460275715007, 1.04673987457171, 1.041263130659057, 1.039329624370229)
$data[16] = @(1.02, 1.040455891908478, 1.043962377020215, 1.038542272542003, 1.037107319114064, 1.042656549162241, 1.046742056550455, 1.047371750483498, 1.041970749007942, 1.040540889892405)
$data[17] = @(1.02, 1.041041070189849, 1.044430252893771, 1.039058070460103, 1.037940146444069, 1.042880123548029, 1.047189007850369, 1.047766974499816, 1.042413287871234, 1.041299239228587)
$data[18] = @(1.02, 1.041381995165399, 1.044702816552199, 1.039358511879345, 1.038425443706658, 1.043010171900553, 1.047449269056848, 1.047997096401108, 1.042670937541212, 1.041741058524147)
$data[19] = @(1.02, 1.041498174448038, 1.044795696316228, 1.03946088469039, 1.038590837458625, 1.043054454375422, 1.04753793749855, 1.048075493554177, 1.04275870910367, 1.041891621205674)
$data[20] = @(1.02, 1.040978327542355, 1.044380089534074, 1.03900277322479, 1.037850841466935, 1.042856173294768, 1.047141099592598, 1.047724612732405, 1.042365856940994, 1.041217928765522)
$data[21] = @(1.02, 1.039283523891498, 1.043024885004104, 1.037508503196552, 1.035439378369131, 1.042207309712374, 1.045845765994136, 1.046579073991724, 1.041083036303017, 1.039021604154995)
$data[22] = @(1.02, 1.038214228281512, 1.042169668734791, 1.036565166579317, 1.033918715264841, 1.04179606345203, 1.045027306806335, 1.04585509975833, 1.040272102514071, 1.037635895793615)
$data[23] = @(1.02, 1.038781438468125, 1.042623337042273, 1.037065614723797, 1.034725281276494, 1.042014386099495, 1.045461573864114, 1.046239249269266, 1.040702412686788, 1.038370948529392)
$data[24] = @(1.02, 1.041006679478637, 1.044402757248873, 1.039027760946687, 1.037891196014702, 1.042866996491891, 1.047162748625099, 1.047743755444759, 1.042387290406636, 1.041254671048468)
$data[25] = @(1.02, 1.043573324109701, 1.046454360300126, 1.041288497790537, 1.041546450985279, 1.043842351847973, 1.049119731292235, 1.049473778155196, 1.044323877600715, 1.044581029646666)

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $rowVals[$i]
    }
}

Write-Host "Updated vm_pu values for rows 2-25 (380 kV slack voltage case)"
